# Adds support for date fields: two new shared-string labelled rows
# ("value" / "budget") at the bottom of the "details" sheet, with the
# sheet's used range/selection/column-width bookkeeping following along
# the way Excel would update them.

$wb = $excel.ActiveWorkbook

$toto    = $wb.Worksheets.Item(1)   # "toto"
$tata    = $wb.Worksheets.Item(2)   # "tata"
$details = $wb.Worksheets.Item(3)   # "details"

# --- new data rows on the "details" sheet ------------------------------
$details.Range("D17").Value = "value"
$details.Range("E17").Value = 50000

$details.Range("D18").Value = "budget"
$details.Range("E18").Value = 300

# --- selection bookkeeping ----------------------------------------------
# Keep the other two sheets' remembered selection where it already was.
$toto.Range("B7").Select()
$tata.Range("A11").Select()

# The active sheet's selection now spans the newly added block.
$details.Range("D17:E18").Select()
$details.Activate()

# --- minor column-width bookkeeping (matches the source file's slight
# global column width nudge that came along with this edit) -------------
$toto.Columns.Item(1).ColumnWidth = 7.5
$tata.Columns.Item(1).ColumnWidth = 7.5
$details.Columns.Item(1).ColumnWidth = 7.666666666666667
